# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> bound to the (only) slide master -> the theme
#                            that actually paints every slide ("Integral").
#   ppt/theme/theme2.xml -> bound to the notes master ("Office Theme").
#
# The authored edit swaps the two themes' contents: the slide master's
# theme becomes the default "Office Theme" palette, and the notes
# master's theme becomes "Integral". Apply the reachable half of that
# swap through the PowerPoint object model: push the "Office Theme"
# color scheme onto the presentation's active theme (theme1.xml) via
# the per-color RGB setters on ThemeColorScheme -- the only writable
# surface this host exposes for theme color data. (dk1/lt1 are already
# black/white in both themes, so only the remaining ten scheme colors
# actually need to move.)

# This host's PowerShell dialect has no built-in RGB() function (that's
# a VBA-ism) -- define the standard r + g*256 + b*65536 packing used by
# ColorFormat.RGB ourselves.
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$colors = $s.ThemeColorScheme

$colors.Colors(1).RGB  = RGB 0 0 0          # dk1     000000
$colors.Colors(2).RGB  = RGB 255 255 255    # lt1     FFFFFF
$colors.Colors(3).RGB  = RGB 68 84 106      # dk2     44546A
$colors.Colors(4).RGB  = RGB 231 230 230    # lt2     E7E6E6
$colors.Colors(5).RGB  = RGB 91 155 213     # accent1 5B9BD5
$colors.Colors(6).RGB  = RGB 237 125 49     # accent2 ED7D31
$colors.Colors(7).RGB  = RGB 165 165 165    # accent3 A5A5A5
$colors.Colors(8).RGB  = RGB 255 192 0      # accent4 FFC000
$colors.Colors(9).RGB  = RGB 68 114 196     # accent5 4472C4
$colors.Colors(10).RGB = RGB 112 173 71     # accent6 70AD47
$colors.Colors(11).RGB = RGB 5 99 193       # hlink   0563C1
$colors.Colors(12).RGB = RGB 149 79 114     # folHlink 954F72
